# Repull data, push all data, mean calculation
# Update column F ("dSF") values for several rows to reflect the repulled data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    3  = 0
    6  = -8
    8  = -6
    13 = -3
    17 = 0
    18 = -3
    19 = -2
    24 = -3
    27 = 1
    28 = 6
    30 = 2
    31 = -1
    37 = -2
    46 = 7
}

foreach ($row in $updates.Keys) {
    $ws.Range("F$row").Value = $updates[$row]
}
